$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Latitude/Longitude columns (H:I) previously held placeholder text
# values (shared, re-used city coordinates). Replace them with the real,
# per-person numeric Latitude/Longitude (used later for Country /
# Country_code lookups in PublicPerson).

$ws.Range("H2").Value = 40.712800000000001
$ws.Range("I2").Value = -74.006

$ws.Range("H3").Value = 34.052199999999999
$ws.Range("I3").Value = -118.2437

$ws.Range("H4").Value = 48.8566
$ws.Range("I4").Value = 2.3521999999999998

$ws.Range("H5").Value = 43.296500000000002
$ws.Range("I5").Value = 5.3697999999999997

$ws.Range("H6").Value = 40.416800000000002
$ws.Range("I6").Value = -3.7038000000000002

$ws.Range("H7").Value = 41.385100000000001
$ws.Range("I7").Value = 2.1734

$ws.Range("H8").Value = -23.5505
$ws.Range("I8").Value = -46.633299999999998

$ws.Range("H9").Value = -22.9068
$ws.Range("I9").Value = -43.172899999999998

# Display the new numeric Latitude/Longitude with six decimal places,
# grouped by each row's original formatting so existing fonts are kept.
$ws.Range("H2:I2").NumberFormat = "0.000000"
$ws.Range("H6:I6").NumberFormat = "0.000000"
$ws.Range("H3:I5").NumberFormat = "0.000000"
$ws.Range("H7:I9").NumberFormat = "0.000000"

# Widen column I to fit the new numeric content, and leave the edited
# range selected.
$ws.Columns.Item(9).AutoFit()
$ws.Range("H2:I9").Select()
